$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F26").Value = 107
$ws.Range("G26").Value = 2740.27
$ws.Range("F29").Value = 46
$ws.Range("G29").Value = 4712.7
$ws.Range("F34").Value = 69
$ws.Range("G34").Value = 2121.06
$ws.Range("B41").Value = 81899.25
$ws.Range("F79").Value = 5
$ws.Range("G79").Value = 1777.05
$ws.Range("B102").Value = 133992.66
$ws.Range("B151").Value = 64196
$ws.Range("F151").Value = 1
$ws.Range("G151").Value = 32143.58
$ws.Range("B152").Value = 65258
$ws.Range("F152").Value = 2
$ws.Range("G152").Value = 64287.16
$ws.Range("F224").Value = 97
$ws.Range("G224").Value = 7625.17
$ws.Range("F231").Value = 23
$ws.Range("G231").Value = 784.53
$ws.Range("B236").Value = 63255
$ws.Range("F236").Value = 92
$ws.Range("G236").Value = 7544
$ws.Range("B237").Value = 57004
$ws.Range("F237").Value = 5
$ws.Range("G237").Value = 410
$ws.Range("B241").Value = 64329
$ws.Range("E241").Value = 128.32
$ws.Range("F241").Value = 1
$ws.Range("G241").Value = 120.69
$ws.Range("B242").Value = 57552
$ws.Range("E242").Value = 136.86
$ws.Range("F242").Value = -5
$ws.Range("G242").Value = -603.45
$ws.Range("B250").Value = 101317.57
$ws.Range("F275").Value = 168
$ws.Range("G275").Value = 9140.879999999999
$ws.Range("F279").Value = 152
$ws.Range("G279").Value = 16029.92
$ws.Range("B283").Value = 114338.32
$ws.Range("F311").Value = 10
$ws.Range("G311").Value = 6765.2
$ws.Range("B318").Value = 24230.03
$ws.Range("B325").Value = 66188
$ws.Range("C325").Value = "HIM-Baby Care Gift Pack (Ww)1"
$ws.Range("D325").Value = 315.8
$ws.Range("E325").Value = 377.31
$ws.Range("F325").Value = 49
$ws.Range("G325").Value = 15474.2
$ws.Range("B326").Value = 48719
$ws.Range("C326").Value = "HIM-BABY CARE GIFT PACK (WW)1"
$ws.Range("D326").Value = 295.75
$ws.Range("E326").Value = 353.35
$ws.Range("F326").Value = -82
$ws.Range("G326").Value = -24251.5
$ws.Range("B370").Value = 64983
$ws.Range("C370").Value = "HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S"
$ws.Range("F370").Value = 6
$ws.Range("G370").Value = 514.08
$ws.Range("B371").Value = 66194
$ws.Range("C371").Value = "HIM-Total Care Baby Pants Diapers-M-9s"
$ws.Range("F371").Value = 39
$ws.Range("G371").Value = 3341.52
$ws.Range("B388").Value = 63531
$ws.Range("E388").Value = 152.53
$ws.Range("F388").Value = 38
$ws.Range("G388").Value = 5452.24
$ws.Range("B389").Value = 57802
$ws.Range("E389").Value = 162.71
$ws.Range("F389").Value = -79
$ws.Range("G389").Value = -11334.92
$ws.Range("F408").Value = 207
$ws.Range("G408").Value = 35465.31
$ws.Range("B412").Value = 51078.17
$ws.Range("B483").Value = 58047
$ws.Range("D483").Value = 105.54
$ws.Range("E483").Value = 126.1
$ws.Range("F483").Value = 34
$ws.Range("G483").Value = 3588.36
$ws.Range("B484").Value = 47097
$ws.Range("D484").Value = 112.28
$ws.Range("E484").Value = 134.16
$ws.Range("F484").Value = 15
$ws.Range("G484").Value = 1684.2
$ws.Range("B553").Value = 65066
$ws.Range("E553").Value = 13.61
$ws.Range("F553").Value = 90
$ws.Range("G553").Value = 1152.9
$ws.Range("B554").Value = 53263
$ws.Range("E554").Value = 15.29
$ws.Range("F554").Value = -309
$ws.Range("G554").Value = -3958.29
$ws.Range("B559").Value = 45706
$ws.Range("E559").Value = 23.58
$ws.Range("F559").Value = -202
$ws.Range("G559").Value = -3985.46
$ws.Range("B560").Value = 64922
$ws.Range("E560").Value = 20.98
$ws.Range("F560").Value = 67
$ws.Range("G560").Value = 1321.91
$ws.Range("B567").Value = 64925
$ws.Range("E567").Value = 13.97
$ws.Range("F567").Value = 111
$ws.Range("G567").Value = 1459.65
$ws.Range("B568").Value = 45709
$ws.Range("E568").Value = 15.69
$ws.Range("F568").Value = -300
$ws.Range("G568").Value = -3945
$ws.Range("F581").Value = 153
$ws.Range("G581").Value = 7406.73
$ws.Range("B587").Value = 42710.56
$ws.Range("F607").Value = 47
$ws.Range("G607").Value = 6639.69
$ws.Range("F612").Value = 83
$ws.Range("G612").Value = 18445.09
$ws.Range("B615").Value = 149076.82
$ws.Range("B672").Value = 64830
$ws.Range("E672").Value = 34.9
$ws.Range("F672").Value = 91
$ws.Range("G672").Value = 2987.53
$ws.Range("B673").Value = 60022
$ws.Range("E673").Value = 37.22
$ws.Range("F673").Value = -113
$ws.Range("G673").Value = -3709.79
$ws.Range("F701").Value = 397
$ws.Range("G701").Value = 6296.42
$ws.Range("F703").Value = 68
$ws.Range("G703").Value = 2936.24
$ws.Range("F705").Value = 64
$ws.Range("G705").Value = 2763.52
$ws.Range("F706").Value = 57
$ws.Range("G706").Value = 1714.56
$ws.Range("F707").Value = 219
$ws.Range("G707").Value = 9583.440000000001
$ws.Range("B708").Value = 41074.25
$ws.Range("F748").Value = 14
$ws.Range("G748").Value = 2098.32
$ws.Range("B755").Value = 77947.35000000001
$ws.Range("F813").Value = 28
$ws.Range("G813").Value = 1255.8
$ws.Range("B829").Value = 65079
$ws.Range("F829").Value = 6
$ws.Range("G829").Value = 245.22
$ws.Range("B830").Value = 65362
$ws.Range("F830").Value = 0
$ws.Range("G830").Value = 0
$ws.Range("B838").Value = 330842.3
$ws.Range("F890").Value = 228
$ws.Range("G890").Value = 6892.44
$ws.Range("F893").Value = 89
$ws.Range("G893").Value = 25175.43
$ws.Range("B897").Value = 339552.49
$ws.Range("F924").Value = 9
$ws.Range("G924").Value = 3642.93
$ws.Range("F934").Value = 30
$ws.Range("G934").Value = 16015.5
$ws.Range("B936").Value = 116639.47
$ws.Range("B942").Value = 5090358.03
$ws.Range("B943").Value = 5090358.03
